$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-metadata")

# Row 7: "Start Data Row" -> Value becomes numeric 4, Example becomes 4
$ws.Range("B7").Value = 4
$ws.Range("D7").Value = 4

# Row 8: "Start Data Col" -> Value becomes "C"
$ws.Range("B8").Value = "C"
